# Apply the Nov 9 2024 crypto-price refresh to the "cryptos" sheet.
# Source data keeps Price/Volume columns as literal text (e.g. "196.36",
# "1.00"), never as numbers, so numeric-looking prices are written through
# Set-TextCell, which forces Text format just long enough to assign the
# string, then restores the cell's original Style so no formatting
# actually changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

# Row 2: Bitcoin
$ws.Cells.Item(2, 4).Value = '75.975.14'
$ws.Cells.Item(2, 5).Value = '  +0.06%  '

# Row 3: Ethereum
$ws.Cells.Item(3, 4).Value = '3.013.40'
$ws.Cells.Item(3, 5).Value = '  +3.89%  '

# Row 4: TetherUSD
$ws.Cells.Item(4, 5).Value = '  +0.01%  '

# Row 5: Solana
Set-TextCell 5 4 '196.36'
$ws.Cells.Item(5, 5).Value = '  -0.29%  '

# Row 6: BNB
Set-TextCell 6 4 '615.16'
$ws.Cells.Item(6, 5).Value = '  +4.34%  '

# Row 7: USDC
$ws.Cells.Item(7, 5).Value = '  +0.02%  '

# Row 8: XRP
Set-TextCell 8 4 '0.548'
$ws.Cells.Item(8, 5).Value = '  +0.57%  '

# Row 9: Dogecoin
Set-TextCell 9 4 '0.204'
$ws.Cells.Item(9, 5).Value = '  +6.36%  '

# Row 10: LidoStakedEther
$ws.Cells.Item(10, 4).Value = '3.016.28'
$ws.Cells.Item(10, 5).Value = '  +3.72%  '

# Row 11: Cardano
$ws.Cells.Item(11, 5).Value = '  -0.33%  '

# Row 12: TRON
$ws.Cells.Item(12, 5).Value = '  -0.26%  '

# Row 13: Toncoin
Set-TextCell 13 4 '5.22'
$ws.Cells.Item(13, 5).Value = '  +7.50%  '

# Row 14: WrappedliquidstakedEther2.0
$ws.Cells.Item(14, 4).Value = '3.583.81'
$ws.Cells.Item(14, 5).Value = '  +3.99%  '

# Row 15: Avalanche
Set-TextCell 15 4 '28.76'
$ws.Cells.Item(15, 5).Value = '  +4.17%  '

# Row 16: WrappedBTC
$ws.Cells.Item(16, 4).Value = '75.980.17'
$ws.Cells.Item(16, 5).Value = '  +0.18%  '

# Row 17: ShibaInu
$ws.Cells.Item(17, 5).Value = '  +2.85%  '

# Row 18: WrappedEther
$ws.Cells.Item(18, 4).Value = '3.020.50'
$ws.Cells.Item(18, 5).Value = '  +4.70%  '

# Row 19: Chainlink
Set-TextCell 19 4 '13.41'
$ws.Cells.Item(19, 5).Value = '  +2.75%  '

# Row 20: Uniswap
Set-TextCell 20 4 '8.88'
$ws.Cells.Item(20, 5).Value = '  +3.00%  '

# Row 21: BitcoinCash
Set-TextCell 21 4 '378.14'
$ws.Cells.Item(21, 5).Value = '  +3.41%  '

# Row 22: SuiNetwork
Set-TextCell 22 4 '2.35'
$ws.Cells.Item(22, 5).Value = '  +6.11%  '

# Row 23: Polkadot
Set-TextCell 23 4 '4.37'
$ws.Cells.Item(23, 5).Value = '  +2.62%  '

# Row 24: WrappedeETH
$ws.Cells.Item(24, 2).Value = 'WrappedeETH'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Cells.Item(24, 4).Value = '3.172.07'
$ws.Cells.Item(24, 5).Value = '  +3.93%  '

# Row 25: Litecoin
$ws.Cells.Item(25, 2).Value = 'Litecoin'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell 25 4 '72.32'
$ws.Cells.Item(25, 5).Value = '  +1.01%  '

# Row 26: Dai
$ws.Cells.Item(26, 5).Value = '  -0.02%  '

# Row 27: NEARProtocol
Set-TextCell 27 4 '4.29'
$ws.Cells.Item(27, 5).Value = '  +2.67%  '

# Row 28: Aptos
Set-TextCell 28 4 '9.72'
$ws.Cells.Item(28, 5).Value = '  +2.57%  '

# Row 29: PEPE
$ws.Cells.Item(29, 5).Value = '  +2.67%  '

# Row 30: Binance-PegBSC-USD
Set-TextCell 30 4 '1.00'
$ws.Cells.Item(30, 5).Value = '  +0.01%  '

# Row 31: InternetComputer(DFINITY)
Set-TextCell 31 4 '8.19'
$ws.Cells.Item(31, 5).Value = '  +2.64%  '

# Row 32: Fetch.AI
Set-TextCell 32 4 '1.38'
$ws.Cells.Item(32, 5).Value = '  +2.19%  '

# Row 33: Bittensor
Set-TextCell 33 4 '489.85'
$ws.Cells.Item(33, 5).Value = '  +0.06%  '

# Row 34: PancakeSwap
$ws.Cells.Item(34, 5).Value = '  +6.38%  '

# Row 35: FirstDigitalUSD
$ws.Cells.Item(35, 5).Value = '  +0.08%  '

# Row 36: Kaspa
$ws.Cells.Item(36, 2).Value = 'Kaspa'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell 36 4 '0.122'
$ws.Cells.Item(36, 5).Value = '  +13.97%  '

# Row 37: EthereumClassic
$ws.Cells.Item(37, 2).Value = 'EthereumClassic'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell 37 4 '20.48'
$ws.Cells.Item(37, 5).Value = '  +2.99%  '

# Row 38: Monero
Set-TextCell 38 4 '161.96'
$ws.Cells.Item(38, 5).Value = '  -2.19%  '

# Row 39: WhiteBITCoin
$ws.Cells.Item(39, 5).Value = '  +1.67%  '

# Row 40: Aave
Set-TextCell 40 4 '189.69'
$ws.Cells.Item(40, 5).Value = '  +6.38%  '

# Row 41: PolygonEcosystemToken
Set-TextCell 41 4 '0.375'
$ws.Cells.Item(41, 5).Value = '  -1.14%  '

# Row 42: Cronos
$ws.Cells.Item(42, 5).Value = '  -6.01%  '

# Row 43: USDe
$ws.Cells.Item(43, 5).Value = '  +0.00%  '

# Row 44: RenderToken
Set-TextCell 44 4 '5.08'
$ws.Cells.Item(44, 5).Value = '  +5.58%  '

# Row 45: Mantle
Set-TextCell 45 4 '0.766'
$ws.Cells.Item(45, 5).Value = '  +18.20%  '

# Row 46: OKB
$ws.Cells.Item(46, 2).Value = 'OKB'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell 46 4 '41.34'
$ws.Cells.Item(46, 5).Value = '  +3.06%  '

# Row 47: ImmutableX
$ws.Cells.Item(47, 2).Value = 'ImmutableX'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 47 4 '1.24'
$ws.Cells.Item(47, 5).Value = '  +6.13%  '

# Row 48: Stacks
$ws.Cells.Item(48, 5).Value = '  +0.48%  '

# Row 49: dogwifhat
Set-TextCell 49 4 '2.40'
$ws.Cells.Item(49, 5).Value = '  +7.96%  '

# Row 50: ARBITRUM
Set-TextCell 50 4 '0.589'
$ws.Cells.Item(50, 5).Value = '  +3.05%  '

# Row 51: Filecoin
Set-TextCell 51 4 '3.84'
$ws.Cells.Item(51, 5).Value = '  +1.11%  '
